# Apply "new satellogic data 20240713" edits to the consolidated sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix / relabel the column headers in row 1 -------------------------
$ws.Range("C1").Value  = "Restricted cash"
$ws.Range("F1").Value  = "Total Current Assets"
$ws.Range("G1").Value  = "Property and equipment,net"
$ws.Range("I1").Value  = "Other non-current assets"
$ws.Range("J1").Value  = "Total Assets"
$ws.Range("O1").Value  = "Operating lease liabilities (current liabilities)"
$ws.Range("P1").Value  = "Contract liabilities (current liabilities)"
$ws.Range("R1").Value  = "Total Current Liabilities"
$ws.Range("T1").Value  = "Contract liabilities"
$ws.Range("V1").Value  = "Total Liabilities"
$ws.Range("AA1").Value = "Total Equity"
$ws.Range("AB1").Value = "Liabilities and Equity"
$ws.Range("AJ1").Value = "Total operating expenses"
$ws.Range("AT1").Value = "Net loss"

# --- 2. Freeze all the interpolation formulas into their cached values ----
$rng = $ws.Range("A1:AT18")
$rng.Value = $rng.Value()

# --- 3. Restore the sheet/window view state --------------------------------
$ws.Range("Q6").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 38   # "AL" is the 38th column -> topLeftCell="AL1"
$win.ScrollRow = 1

$wb.Windows.Item(1).WindowState = $wb.Windows.Item(1).WindowState
